$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, shifting existing rows 37-40 down to 38-41
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).RowHeight = 13.05

# Populate the newly inserted row 37
$ws.Range("A37").Value = "Olivers"
$ws.Range("B37").Value = "Pietrs, Josh"
$ws.Range("C37").Value = "013"
$ws.Range("E37").Value = "0008284"
$ws.Range("F37").Font.Name = "Arial"
